$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 335; this shifts rows 335..453 down to 336..454
# and the sheet dimension grows from A1:R453 to A1:R454, matching the diff.
$ws.Rows(335).Insert()

# Populate the newly inserted row 335 with the "Ajo" record added by this edit.
# (Same attributes as the record that was at row 335 before the insert, except
# the date (D) and volume (J) values, which are new.)
$ws.Range("A335").Value = 8
$ws.Range("B335").Value = "Terminal La Palmera de La Serena"
$ws.Range("C335").Value = "Coquimbo"
$ws.Range("D335").Value = 45027
$ws.Range("E335").Value = 4
$ws.Range("F335").Value = 100112003
$ws.Range("G335").Value = "Ajo"
$ws.Range("H335").Value = "Chino"
$ws.Range("I335").Value = "Primera"
$ws.Range("J335").Value = 400
$ws.Range("K335").Value = 17000
$ws.Range("L335").Value = 18000
$ws.Range("M335").Value = 17500
$ws.Range("N335").Value = "$/caja 10 kilos"
$ws.Range("O335").Value = "China"
$ws.Range("P335").Value = 1750
$ws.Range("Q335").Value = 10
$ws.Range("R335").Value = "Hortaliza"
